# Cleaning on merged file
# Applies the data-cleaning edits made to the "Cleaning" sheet:
#   1. Seven "####" roll-number cells in column AL were stored as text with a
#      stray "#10" suffix (e.g. "1028#10"); clean them into plain numbers.
#   2. Seven rows were missing their column R ("Total III") value; fill them in.
#   3. Row 59 had a copy/paste slip: BY59 (the "re-exam"/back-log marks column)
#      was left at 0 and BZ59 ("Total IV") was a stale literal instead of the
#      live SUM of BX59:BY59; fix both.
#   4. Add a new helper column CC "Company -1/0" that flags (1/0) whether a
#      candidate has a placement Company recorded in column CB.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Strip the "#10" suffix / fix the 7 text-as-number cells in column AL ---
$ws.Range("AL17").Value = 1028
$ws.Range("AL24").Value = 1140
$ws.Range("AL27").Value = 512
$ws.Range("AL29").Value = 506
$ws.Range("AL32").Value = 1048
$ws.Range("AL35").Value = 1094
$ws.Range("AL43").Value = 496

# --- 2. Fill the 7 missing "Total III" (column R) values ---
$ws.Range("R10").Value = 419
$ws.Range("R13").Value = 394
$ws.Range("R23").Value = 427
$ws.Range("R57").Value = 416
$ws.Range("R67").Value = 420
$ws.Range("R73").Value = 437
$ws.Range("R80").Value = 423

# --- 3. Row 59 fix: real BY59 value + live total formula in BZ59 ---
$ws.Range("BY59").Value = 409
$ws.Range("BZ59").Formula = "=SUM(BX59,BY59)"

# --- 4. New column CC: "Company -1/0" flag ---
$ws.Range("CC1").Value = "Company -1/0"
$ws.Range("CC1").Font.Bold = $true
$ws.Range("CC1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("CC1").VerticalAlignment = -4160     # xlTop
$ws.Range("CC1").Borders.Item(7).LineStyle = 1   # xlEdgeLeft  / xlContinuous
$ws.Range("CC1").Borders.Item(10).LineStyle = 1  # xlEdgeRight / xlContinuous

# Row-by-row flag: 1 when the candidate's row has a Company in column CB, else 0.
$ccFlags = "0,1,0,1,0,1,1,0,1,0,0,0,0,1,0,1,0,1,1,0,0,1,1,0,0,0,1,1,1,0,0,1,1,1,1,1,0,1,1,1,1,1,0,1,1,0,1,1,0,0,1,0,0,1,0,1,0,0,1,0,0,1,0,0,0,0,0,1,1,0,1,1,1,0,0,0,0,0,0,0,1,1,0,0,0,1,1,1,0,0,0,1,1,0,1,1,0,1,1,1,1,1,0,0,1,0,1,0,0,0,0,1,0,1,0"
$parts = $ccFlags.Split(",")
$n = $parts.Length
$ccVals = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $ccVals[$i,0] = [int]$parts[$i]
}
$ws.Range("CC2:CC116").Value = $ccVals
